$d = $word.ActiveDocument

# 1) Fix the title: merge the split "ECM1410 Cover " / "page" / " " runs
#    (with gramStart/gramEnd proofErr markers) into a single run of text.
$d.Content.Find.Execute("ECM1410 Cover page ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ECM1410 Cover page ", 2) | Out-Null

# 2) Second table (the schedule table) updates.
$tbl = $d.Tables.Item(2)

# Row 3 (22/3/23 -> 22/2/23)
$tbl.Cell(3, 1).Range.Text = "22/2/23"

# Row 4 (previously blank) -> fill in meeting details
$tbl.Cell(4, 1).Range.Text = "24/2/23"
$tbl.Cell(4, 2).Range.Text = "11:40"
$tbl.Cell(4, 3).Range.Text = "1h15"
$tbl.Cell(4, 4).Range.Text = "Driver"
$tbl.Cell(4, 5).Range.Text = "Observer"
